# Scheduled-runner data refresh: updates cached market-board price columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) on the leve
# worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with freshly pulled
# values. These are static snapshot values (no formulas backing them), so
# each touched cell is simply overwritten with its new value.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2163.1135
$ws.Range("I15").Value = 2163.1135
$ws.Range("K15").Value = 6489.3405
$ws.Range("M15").Value = -6320.3405

$ws.Range("H34").Value = 12537.714
$ws.Range("I34").Value = 7960.8335
$ws.Range("K34").Value = 7960.8335
$ws.Range("M34").Value = -7757.8335

$ws.Range("H36").Value = 12537.714
$ws.Range("I36").Value = 7960.8335
$ws.Range("K36").Value = 7960.8335
$ws.Range("M36").Value = -7245.8335

$ws.Range("H64").Value = 7331.6665
$ws.Range("J64").Value = 7331.6665
$ws.Range("L64").Value = 7331.6665
$ws.Range("N64").Value = -7827.6665

$ws.Range("H67").Value = 7331.6665
$ws.Range("J67").Value = 7331.6665
$ws.Range("L67").Value = 7331.6665
$ws.Range("N67").Value = -9047.666499999999

$ws.Range("H138").Value = 3260.5151
$ws.Range("I138").Value = 2585.963
$ws.Range("J138").Value = 6296
$ws.Range("K138").Value = 7757.889000000001
$ws.Range("L138").Value = 18888
$ws.Range("M138").Value = -2617.889000000001
$ws.Range("N138").Value = -29168

$ws.Range("H141").Value = 7155
$ws.Range("I141").Value = 6566.121
$ws.Range("K141").Value = 19698.363
$ws.Range("M141").Value = -14518.363

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 700
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 700
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 700
$ws.Range("N4").Value = -932
$ws.Range("M4").ClearContents()

$ws.Range("H5").Value = 466.55554
$ws.Range("I5").Value = 529.8
$ws.Range("K5").Value = 529.8
$ws.Range("M5").Value = -417.8

$ws.Range("H6").Value = 4875
$ws.Range("J6").Value = 9500
$ws.Range("L6").Value = 9500
$ws.Range("N6").Value = -9846

$ws.Range("H32").Value = 10703.765
$ws.Range("I32").Value = 7824.269
$ws.Range("J32").Value = 20062.125
$ws.Range("K32").Value = 7824.269
$ws.Range("L32").Value = 20062.125
$ws.Range("M32").Value = -7537.269
$ws.Range("N32").Value = -20636.125

$ws.Range("H39").Value = 15000.125
$ws.Range("I39").Value = 15166.667
$ws.Range("K39").Value = 15166.667
$ws.Range("M39").Value = -14646.667

$ws.Range("H45").Value = 361250.3
$ws.Range("I45").Value = 1252649
$ws.Range("J45").Value = 4690.85
$ws.Range("K45").Value = 1252649
$ws.Range("L45").Value = 4690.85
$ws.Range("M45").Value = -1252272
$ws.Range("N45").Value = -5444.85

$ws.Range("H61").Value = 3957.5732
$ws.Range("I61").Value = 3779.5857
$ws.Range("J61").Value = 4995.8335
$ws.Range("K61").Value = 3779.5857
$ws.Range("L61").Value = 4995.8335
$ws.Range("M61").Value = -3567.5857
$ws.Range("N61").Value = -5419.8335

$ws.Range("H110").Value = 4676.4614
$ws.Range("I110").Value = 4143.778
$ws.Range("K110").Value = 4143.778
$ws.Range("M110").Value = -2098.778

$ws.Range("H134").Value = 629110.25
$ws.Range("J134").Value = 447930.5
$ws.Range("L134").Value = 447930.5
$ws.Range("N134").Value = -458070.5

$ws.Range("H136").Value = 3957.5732
$ws.Range("I136").Value = 3779.5857
$ws.Range("J136").Value = 4995.8335
$ws.Range("K136").Value = 11338.7571
$ws.Range("L136").Value = 14987.5005
$ws.Range("M136").Value = -8788.757100000001
$ws.Range("N136").Value = -20087.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 466.55554
$ws.Range("I4").Value = 529.8
$ws.Range("K4").Value = 529.8
$ws.Range("M4").Value = -414.8

$ws.Range("H21").Value = 27654.223
$ws.Range("J21").Value = 27654.223
$ws.Range("L21").Value = 27654.223
$ws.Range("N21").Value = -28126.223

$ws.Range("H105").Value = 2570.4285
$ws.Range("I105").Value = 1999
$ws.Range("K105").Value = 1999
$ws.Range("M105").Value = -252

$ws.Range("H134").Value = 2487.6511
$ws.Range("I134").Value = 2334.475
$ws.Range("K134").Value = 7003.424999999999
$ws.Range("M134").Value = -4468.424999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6229.35
$ws.Range("I31").Value = 4535
$ws.Range("J31").Value = 7615.636
$ws.Range("K31").Value = 4535
$ws.Range("L31").Value = 7615.636
$ws.Range("M31").Value = -4240
$ws.Range("N31").Value = -8205.636

$ws.Range("H34").Value = 6229.35
$ws.Range("I34").Value = 4535
$ws.Range("J34").Value = 7615.636
$ws.Range("K34").Value = 4535
$ws.Range("L34").Value = 7615.636
$ws.Range("M34").Value = -4333
$ws.Range("N34").Value = -8019.636

$ws.Range("H58").Value = 65461.062
$ws.Range("I58").Value = 74468.78999999999
$ws.Range("J58").Value = 2407
$ws.Range("K58").Value = 74468.78999999999
$ws.Range("L58").Value = 2407
$ws.Range("M58").Value = -74265.78999999999
$ws.Range("N58").Value = -2813

$ws.Range("H132").Value = 2631
$ws.Range("I132").Value = 2785.6667
$ws.Range("J132").Value = 2074.2
$ws.Range("K132").Value = 8357.000100000001
$ws.Range("L132").Value = 6222.599999999999
$ws.Range("M132").Value = -5827.000100000001
$ws.Range("N132").Value = -11282.6

$ws.Range("H136").Value = 65461.062
$ws.Range("I136").Value = 74468.78999999999
$ws.Range("J136").Value = 2407
$ws.Range("K136").Value = 223406.37
$ws.Range("L136").Value = 7221
$ws.Range("M136").Value = -220856.37
$ws.Range("N136").Value = -12321

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 7359.6
$ws.Range("I58").Value = 6449.75
$ws.Range("K58").Value = 19349.25
$ws.Range("M58").Value = -19221.25

$ws.Range("H97").Value = 549.2
$ws.Range("I97").Value = 430.66666
$ws.Range("K97").Value = 1291.99998
$ws.Range("M97").Value = -795.9999800000001

$ws.Range("H132").Value = 2080
$ws.Range("I132").Value = 1600
$ws.Range("K132").Value = 14400
$ws.Range("M132").Value = -11870

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H119").Value = 51666.668
$ws.Range("J119").Value = 51666.668
$ws.Range("L119").Value = 51666.668
$ws.Range("N119").Value = -61342.668

$ws.Range("H132").Value = 25264.043
$ws.Range("I132").Value = 27504.05
$ws.Range("K132").Value = 82512.14999999999
$ws.Range("M132").Value = -79982.14999999999

$ws.Range("H141").Value = 68500
$ws.Range("J141").Value = 68500
$ws.Range("L141").Value = 68500
$ws.Range("N141").Value = -78860

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2389.7693
$ws.Range("I22").Value = 1356.5714
$ws.Range("J22").Value = 3595.1667
$ws.Range("K22").Value = 1356.5714
$ws.Range("L22").Value = 3595.1667
$ws.Range("M22").Value = -1061.5714
$ws.Range("N22").Value = -4185.1667

$ws.Range("H27").Value = 2389.7693
$ws.Range("I27").Value = 1356.5714
$ws.Range("J27").Value = 3595.1667
$ws.Range("K27").Value = 1356.5714
$ws.Range("L27").Value = 3595.1667
$ws.Range("M27").Value = -1249.5714
$ws.Range("N27").Value = -3809.1667

$ws.Range("H46").Value = 17285.715
$ws.Range("I46").Value = 47000.5
$ws.Range("J46").Value = 5399.8
$ws.Range("K46").Value = 47000.5
$ws.Range("L46").Value = 5399.8
$ws.Range("M46").Value = -46812.5
$ws.Range("N46").Value = -5775.8

$ws.Range("H122").Value = 4696.879
$ws.Range("I122").Value = 4285.643
$ws.Range("K122").Value = 12856.929
$ws.Range("M122").Value = -10406.929

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 6666.3335
$ws.Range("I5").Value = 9999
$ws.Range("J5").Value = 5000
$ws.Range("K5").Value = 9999
$ws.Range("L5").Value = 5000
$ws.Range("M5").Value = -9887
$ws.Range("N5").Value = -5224
